$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel coercing
# numeric-looking strings (e.g. "595.34") into the Number type, and
# without leaving a residual number-format style behind (so the cell
# round-trips exactly like the original unstyled inline-string cells).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '65.594.68'
Set-TextValue $ws.Range("E2") '  -0.05%  '
Set-TextValue $ws.Range("D3") '2.648.07'
Set-TextValue $ws.Range("E3") '  -0.84%  '
Set-TextValue $ws.Range("E4") '  +0.08%  '
Set-TextValue $ws.Range("D5") '595.34'
Set-TextValue $ws.Range("E5") '  -0.92%  '
Set-TextValue $ws.Range("D6") '156.13'
Set-TextValue $ws.Range("E6") '  -0.45%  '
Set-TextValue $ws.Range("E7") '  +0.07%  '
Set-TextValue $ws.Range("E8") '  +3.08%  '
Set-TextValue $ws.Range("E9") '  +3.67%  '
Set-TextValue $ws.Range("E10") '  -0.65%  '
Set-TextValue $ws.Range("E11") '  -1.93%  '
Set-TextValue $ws.Range("E12") '  +1.12%  '
Set-TextValue $ws.Range("D13") '28.73'
Set-TextValue $ws.Range("E13") '  -2.31%  '
Set-TextValue $ws.Range("D14") '0.0000196'
Set-TextValue $ws.Range("E14") '  +0.25%  '
Set-TextValue $ws.Range("D15") '3.128.05'
Set-TextValue $ws.Range("E15") '  -0.62%  '
Set-TextValue $ws.Range("D16") '65.461.42'
Set-TextValue $ws.Range("E16") '  +0.11%  '
Set-TextValue $ws.Range("D17") '2.655.41'
Set-TextValue $ws.Range("E17") '  -0.51%  '
Set-TextValue $ws.Range("D18") '12.58'
Set-TextValue $ws.Range("E18") '  +0.67%  '
Set-TextValue $ws.Range("D19") '4.73'
Set-TextValue $ws.Range("E19") '  -1.73%  '
Set-TextValue $ws.Range("D20") '7.44'
Set-TextValue $ws.Range("E20") '  -1.52%  '
Set-TextValue $ws.Range("D21") '348.28'
Set-TextValue $ws.Range("E21") '  -0.74%  '
Set-TextValue $ws.Range("D23") '69.02'
Set-TextValue $ws.Range("E23") '  -1.02%  '
Set-TextValue $ws.Range("D24") '0.0000111'
Set-TextValue $ws.Range("E24") '  +2.01%  '
Set-TextValue $ws.Range("D25") '9.67'
Set-TextValue $ws.Range("E25") '  -0.20%  '
Set-TextValue $ws.Range("E26") '  +1.39%  '
Set-TextValue $ws.Range("E27") '  -0.52%  '
Set-TextValue $ws.Range("E28") '  -2.10%  '
Set-TextValue $ws.Range("E29") '  +0.00%  '
Set-TextValue $ws.Range("D30") '7.87'
Set-TextValue $ws.Range("E30") '  -3.05%  '
Set-TextValue $ws.Range("E31") '  -1.21%  '
Set-TextValue $ws.Range("D32") '528.58'
Set-TextValue $ws.Range("E32") '  -2.55%  '
Set-TextValue $ws.Range("E33") '  -0.47%  '
Set-TextValue $ws.Range("D34") '6.41'
Set-TextValue $ws.Range("E34") '  -1.99%  '
Set-TextValue $ws.Range("D35") '5.41'
Set-TextValue $ws.Range("E35") '  -0.91%  '
Set-TextValue $ws.Range("E36") '  -0.52%  '
Set-TextValue $ws.Range("D37") '20.34'
Set-TextValue $ws.Range("E37") '  -0.37%  '
Set-TextValue $ws.Range("E38") '  +0.00%  '
Set-TextValue $ws.Range("D39") '156.71'
Set-TextValue $ws.Range("E39") '  -1.32%  '
Set-TextValue $ws.Range("E40") '  -1.36%  '
Set-TextValue $ws.Range("D41") '0.999'
Set-TextValue $ws.Range("D42") '160.70'
Set-TextValue $ws.Range("E42") '  -3.05%  '
Set-TextValue $ws.Range("D43") '4.06'
Set-TextValue $ws.Range("E43") '  -0.26%  '
Set-TextValue $ws.Range("D44") '0.0604'
Set-TextValue $ws.Range("E44") '  -1.04%  '
Set-TextValue $ws.Range("D45") '2.26'
Set-TextValue $ws.Range("E45") '  -0.03%  '
Set-TextValue $ws.Range("D46") '22.57'
Set-TextValue $ws.Range("E46") '  -2.12%  '
Set-TextValue $ws.Range("D47") '0.633'
Set-TextValue $ws.Range("E47") '  -2.10%  '
Set-TextValue $ws.Range("E48") '  -2.03%  '
Set-TextValue $ws.Range("D49") '0.0992'
Set-TextValue $ws.Range("E49") '  -1.10%  '
Set-TextValue $ws.Range("D50") '0.0₆0251'
Set-TextValue $ws.Range("E50") '  +10.42%  '
Set-TextValue $ws.Range("D51") '19.70'
Set-TextValue $ws.Range("E51") '  -1.29%  '
